$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two header-only rows that carried no numeric data.
# Delete the lower one first so the upper row's index doesn't shift
# before it is removed.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()

# Row 2's B and F labels were placeholder "unnamed" headers; correct
# them to "total" (matching column C, which already reads "total").
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
